$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Week 2 table (row 11) ---
$ws.Range("C11").Value = 3
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 5

# --- New "Sleep study" label under Week 3 heading (row 14) ---
$rng = $ws.Range("E14")
$rng.Value = "Sleep study"
$rng.Font.Color = 24832
$rng.Interior.Color = 13561798

# --- Week 3 table (row 16) ---
$ws.Range("B16").Value = 3
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 8

# --- Week 4 table (row 21) ---
$ws.Range("A21").Value = 3
$ws.Range("B21").Value = 3
$ws.Range("C21").Value = 3
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 3
# total for week 4 now exceeds the "good" threshold -> flip the Total cell to the "Bad" look
$tot21 = $ws.Range("H21")
$tot21.Font.Color = 393372
$tot21.Interior.Color = 13551615

# --- Week 5 table (row 26) ---
$ws.Range("A26").Value = 2
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 3
# total for week 5 now exceeds the "good" threshold -> flip the Total cell to the "Bad" look
$tot26 = $ws.Range("H26")
$tot26.Font.Color = 393372
$tot26.Interior.Color = 13551615

# --- Week 6 table (row 31) ---
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = 3
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 7

# --- Move the active selection like the author left it ---
$ws.Range("K32").Select
